# Applies the cryptos-list refresh described in the commit:
#   "Updated cryptos list on Tue Jan  9 12:56:45 UTC 2024 with GitHub Actions"
# All Price/Volume cells are stored as TEXT in the sheet (t="inlineStr"), so for
# any new value that LOOKS like a plain number we briefly force the cell to Text
# format before assigning it (otherwise COM Excel would auto-convert it to a real
# number and silently drop things like trailing zeros), then restore the cell
# style to Normal so no formatting changes leak into the saved file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '46.858.45'
$ws.Cells.Item(2, 5).Value = '  +4.21%  '
$ws.Cells.Item(3, 4).Value = '2.295.80'
$ws.Cells.Item(3, 5).Value = '  +1.47%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '303.77'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.86%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '101.41'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +8.15%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.567'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.30%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.523'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.19%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '36.22'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +5.86%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0788'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.16%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.27'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.33%  '
$ws.Cells.Item(13, 5).Value = '  -0.23%  '
$ws.Cells.Item(14, 4).Value = '2.646.27'
$ws.Cells.Item(14, 5).Value = '  +1.75%  '
$ws.Cells.Item(15, 4).Value = '2.296.99'
$ws.Cells.Item(15, 5).Value = '  +1.65%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.81'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.90%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.808'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.32%  '
$ws.Cells.Item(18, 4).Value = '46.825.30'
$ws.Cells.Item(18, 5).Value = '  +4.62%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.98'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.94%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0940'
$ws.Cells.Item(20, 5).Value = '  +2.14%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.98'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -1.61%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '65.82'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.95%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '251.04'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +4.78%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.88'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.10%  '
$ws.Cells.Item(25, 5).Value = '  -0.08%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.91'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.83%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '42.03'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +5.33%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.23'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -2.46%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.84'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +3.11%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.00'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.20%  '
$ws.Cells.Item(31, 5).Value = '  +10.48%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.56'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.38%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '147.56'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -3.74%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0785'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.66%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.24'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +11.34%  '
$ws.Cells.Item(36, 5).Value = '  +9.69%  '
$ws.Cells.Item(37, 5).Value = '  -0.19%  '
$ws.Cells.Item(38, 2).Value = 'Celestia'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.02'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +15.05%  '
$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.75'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +0.77%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.95'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +6.09%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.32'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +1.95%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0300'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -1.31%  '
$ws.Cells.Item(43, 5).Value = '  +0.01%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.96'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +5.92%  '
$ws.Cells.Item(45, 4).Value = '1.816.57'
$ws.Cells.Item(45, 5).Value = '  +1.56%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '90.01'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +19.26%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.193'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.65%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '73.13'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.69%  '
$ws.Cells.Item(49, 5).Value = '  +5.13%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '95.43'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.96%  '
$ws.Cells.Item(51, 4).Value = '2.521.76'
$ws.Cells.Item(51, 5).Value = '  +1.62%  '
